# The script used to write both an Excel column (A) and an Outlook-facing
# column (B) with the same label strings. It now only writes data to the
# Excel sheet, so column B is no longer populated and the "Test Field #n"
# labels in column A are replaced with the raw numeric values they
# described.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B duplicated column A purely for the Outlook step that no longer
# runs - clear it out.
$ws.Range("B1:B8").ClearContents() | Out-Null

# Column A keeps the header/details text as-is...
$ws.Range("A1").Value = "#22222222222222"
$ws.Range("A2").Value = "This is a heading. "
$ws.Range("A4").Value = "More details. Test # 1. "

# ...but the test-field rows now hold the plain numeric values instead of
# the "Test Field #n: <value>" label strings.
$ws.Range("A6").Value = 204.33
$ws.Range("A7").Value = 201.23231
$ws.Range("A8").Value = 701.9

$ws.Range("A1:C8").Select() | Out-Null
